# Applies the Feb 17 2023 symbol-list price/volume refresh (GitHub Actions bot commit).
# Source cells store numeric-looking values as TEXT (e.g. "309.88", "-3.45%"), so each
# write goes through a text NumberFormat to stop Excel auto-converting to Number/Percent,
# then the style is reset back to "Normal" so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '309.88'
Set-TextValue $ws.Range("E2") '-3.45%'
Set-TextValue $ws.Range("D3") '51.00'
Set-TextValue $ws.Range("E3") '4.47%'
Set-TextValue $ws.Range("E4") '-1.48%'
Set-TextValue $ws.Range("D5") '0.07784'
Set-TextValue $ws.Range("E5") '-3.81%'
Set-TextValue $ws.Range("D6") '4.493'
Set-TextValue $ws.Range("E6") '-2.25%'
Set-TextValue $ws.Range("D7") '1.356'
Set-TextValue $ws.Range("E7") '12.26%'
Set-TextValue $ws.Range("D8") '1.565'
Set-TextValue $ws.Range("E8") '-4.88%'
Set-TextValue $ws.Range("D9") '0.1215'
Set-TextValue $ws.Range("E9") '-5.97%'
Set-TextValue $ws.Range("D10") '0.1980'
Set-TextValue $ws.Range("E10") '1.94%'
Set-TextValue $ws.Range("D11") '0.04799'
Set-TextValue $ws.Range("E11") '3.88%'
Set-TextValue $ws.Range("D12") '0.09501'
Set-TextValue $ws.Range("D13") '0.1044'
Set-TextValue $ws.Range("E13") '-0.75%'
Set-TextValue $ws.Range("D14") '0.001258'
Set-TextValue $ws.Range("E14") '-4.98%'
Set-TextValue $ws.Range("D15") '0.005786'
Set-TextValue $ws.Range("E15") '-1.97%'
Set-TextValue $ws.Range("E16") '2,014.46%'
Set-TextValue $ws.Range("D17") '3.330'
Set-TextValue $ws.Range("E17") '-0.20%'
Set-TextValue $ws.Range("D18") '2.438'
Set-TextValue $ws.Range("E18") '0.30%'
Set-TextValue $ws.Range("E19") '2.67%'
Set-TextValue $ws.Range("D20") '8.023'
Set-TextValue $ws.Range("E20") '-0.69%'
Set-TextValue $ws.Range("D21") '0.1370'
Set-TextValue $ws.Range("E21") '-0.82%'
Set-TextValue $ws.Range("D22") '0.3091'
Set-TextValue $ws.Range("E22") '-1.12%'
Set-TextValue $ws.Range("D23") '0.04149'
Set-TextValue $ws.Range("E23") '-0.39%'
Set-TextValue $ws.Range("D24") '0.001267'
Set-TextValue $ws.Range("E24") '-2.92%'
Set-TextValue $ws.Range("D25") '0.003951'
Set-TextValue $ws.Range("E25") '-7.00%'
Set-TextValue $ws.Range("D26") '0.0001348'
Set-TextValue $ws.Range("E26") '-0.21%'
Set-TextValue $ws.Range("D38") '0.02606'
Set-TextValue $ws.Range("E38") '-4.33%'
Set-TextValue $ws.Range("E39") '4.72%'
Set-TextValue $ws.Range("E40") '74.25%'
Set-TextValue $ws.Range("D41") '0.007869'
Set-TextValue $ws.Range("E41") '0.04%'
Set-TextValue $ws.Range("D42") '0.1424'
Set-TextValue $ws.Range("E42") '-1.19%'
Set-TextValue $ws.Range("D43") '0.008349'
Set-TextValue $ws.Range("E43") '8.37%'
Set-TextValue $ws.Range("D44") '0.007663'
Set-TextValue $ws.Range("E44") '-5.41%'
Set-TextValue $ws.Range("D45") '0.3382'
Set-TextValue $ws.Range("E45") '5.91%'
Set-TextValue $ws.Range("D46") '0.00007268'
Set-TextValue $ws.Range("E46") '5.53%'
Set-TextValue $ws.Range("D47") '0.00000000749'
Set-TextValue $ws.Range("E47") '-0.21%'
Set-TextValue $ws.Range("D48") '0.002616'
Set-TextValue $ws.Range("E48") '-34.66%'
Set-TextValue $ws.Range("D49") '0.05318'
Set-TextValue $ws.Range("E49") '-14.11%'
Set-TextValue $ws.Range("D50") '0.00002097'
Set-TextValue $ws.Range("E50") '-0.21%'
Set-TextValue $ws.Range("D51") '0.0001997'
Set-TextValue $ws.Range("E51") '-0.21%'
